$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = 5.792
$ws.Range("B12").Value = 5.315
$ws.Range("B18").Value = 5.154
$ws.Range("B37").Value = 8.73
$ws.Range("B55").Value = 4.684
$ws.Range("B68").Value = 5.356
$ws.Range("B77").Value = 6.114
$ws.Range("B78").Value = 7.811
